# "optimize the ui verify actions"
# - Move the "Wait" flag on the PageModules sheet from the
#   verify_title_contains row (row 5) up to the click row (row 4).
# - Make PageModules the active sheet/tab, with the selection parked on I4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # "PageModules"

# Move the Wait flag: clear H5, set H4 = 1
$ws.Range("H5").Value = $null
$ws.Range("H4").Value = 1

# Make PageModules the active sheet (updates tabSelected + workbook activeTab)
# and leave the selection on I4.
$ws.Activate()
$ws.Range("I4").Select() | Out-Null
